$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(10.90566285846436, 10.33793240923849, 9.973432640156052, 9.821050522892273, 9.795520433261276, 9.971392905749543, 10.71329654532148, 12.03640013990192, 12.92229191781887, 13.30571129230105, 13.44803193334883, 13.41750915132639, 13.3174779400465, 13.2558304322658, 12.89683549565536, 12.67154535437847, 12.54012429702368, 12.49531326476938, 12.69571876867622, 13.34693787158305, 13.75994557307548, 13.53912656975063, 12.68479587784511, 11.69325189652054)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $colB[$i]
}

$colC = @(8.485397513044358, 8.077531742741494, 7.814714950598358, 7.704577871896445, 7.686108794105026, 7.813241788225763, 8.347378290667491, 9.29379177324404, 9.924610504931852, 10.19714061523203, 10.29823784442984, 10.27655871760724, 10.20550025579777, 10.16170014282522, 9.906507017784875, 9.746235481777413, 9.652694372892618, 9.620790838713233, 9.763437347660137, 10.22642914987654, 10.51674795166079, 10.36293031904897, 9.755664739083148, 9.048916188060856)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $colC[$i]
}

$colD = @(5.961199038305452, 5.837777789821213, 5.762429399052676, 5.731881702214228, 5.726820218134253, 5.762016720213923, 5.918582469276722, 6.226954770516667, 6.451641892973459, 6.552917318722347, 6.591091710825332, 6.582878658009104, 6.556061715270646, 6.539611390216783, 6.445000840687581, 6.386690084790536, 6.353065541060382, 6.341667437746837, 6.392906548936415, 6.563943610236711, 6.674680323193234, 6.615686925350658, 6.39009639662535, 6.143672716418735)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $colD[$i]
}

$colE = @(12.03872722669509, 11.92664659270388, 11.86203899902303, 11.83679769505758, 11.83267280444352, 11.86169415040113, 11.99922562282198, 12.30097141811795, 12.54022165859255, 12.65241744004726, 12.69534603785706, 12.68608157108798, 12.65594056883189, 12.63753474114497, 12.53295355130058, 12.46962823547821, 12.43352391035584, 12.42135557093575, 12.476336628661, 12.66478203025822, 12.79050186253409, 12.72318219347735, 12.47330281956942, 12.21610710934766)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $colE[$i]
}

$colG = @(21.92890485486357, 21.74777922227914, 21.6456564099924, 21.60637356215696, 21.59999285905677, 21.64511712297813, 21.86459656305553, 22.36474795624078, 22.77121069250364, 22.9637542322323, 23.03768985969285, 23.02172219776826, 22.96981688599758, 22.93815441023468, 22.75877490529571, 22.65063578347842, 22.58916096238142, 22.56847323714921, 22.66207299041887, 22.9850355657858, 23.20204009357448, 23.08570351996959, 22.65690005606879, 22.22232244255398)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item(2 + $i, 7).Value = $colG[$i]
}

$colH = @(12.26453792122942, 12.2952220396429, 12.31662140449321, 12.32598408913941, 12.32757750869242, 12.31674507415417, 12.27458591357164, 12.21227116628366, 12.17896802626595, 12.16653983256499, 12.16222578125618, 12.16313743195334, 12.16617704739463, 12.1680900042224, 12.17983508074888, 12.18773802791715, 12.19253971651894, 12.19420944631381, 12.18687023122855, 12.16527358711418, 12.15344550091331, 12.15954889450877, 12.18726175794782, 12.22694192143997)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Cells.Item(2 + $i, 8).Value = $colH[$i]
}

$colI = @(16.99738304866474, 17.07954403260053, 17.1342293487641, 17.15757672210787, 17.16151762891384, 17.13453992096733, 17.0248306948332, 16.84345550159587, 16.73096970978148, 16.68434749003045, 16.66735028520333, 16.67098163230081, 16.68293591948914, 16.69034401372181, 16.73410840945464, 16.76212389726811, 16.7786655409211, 16.78433963562544, 16.75909729108381, 16.67940677930117, 16.63115919331245, 16.65655776669305, 16.76046426446885, 16.88888833418602)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $colI[$i]
}

$colM = @(13.9815265574655, 13.70000007311172, 13.52727121263977, 13.45701633373063, 13.44536161577374, 13.52632305535737, 13.88448927960382, 14.58356249975704, 15.08960747070554, 15.31708035853388, 15.40274076324959, 15.38431486707642, 15.32413775943335, 15.28721274110136, 15.07467917799171, 14.94353473099894, 14.86785183639108, 14.84218632211091, 14.95752201049661, 15.34182688590805, 15.59016615057519, 15.45790853236726, 14.95119924628187, 14.39538545210272)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item(2 + $i, 13).Value = $colM[$i]
}

$colN = @(16.22224096298492, 16.2693648347987, 16.30001671691922, 16.3129403557349, 16.31511247787685, 16.30018925634099, 16.23813327866921, 16.13003285211069, 16.05884537415454, 16.02823763644902, 16.01690188316217, 16.01933192086641, 16.02729993548375, 16.03221373004797, 16.06088134149225, 16.07892236596245, 16.08946626897776, 16.09306498664595, 16.07698456989909, 16.02495262794248, 15.99243124907598, 16.00965288969222, 16.0778601120463, 16.15782727690948)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item(2 + $i, 14).Value = $colN[$i]
}

$colO = @(17.79067923743658, 17.80293235561032, 17.81570205255757, 17.82222221743801, 17.82338429612913, 17.815784660012, 17.79381381437125, 17.79244518700564, 17.81694096451936, 17.83361858467035, 17.84072792013531, 17.83916153052431, 17.83418757792732, 17.83124420123491, 17.81596226307938, 17.80800368365898, 17.80394721684963, 17.80266330997098, 17.80879696750719, 17.83562702286745, 17.85778779889656, 17.84553780508072, 17.80843670698334, 17.78833868404448)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $colO[$i]
}

